$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "Wh"
$ws.Range("D3").Value = 0.95
$ws.Range("E3").Value = 4500
$ws.Range("F3").Value = "W"
$ws.Range("G3").Value = 0.95
$ws.Range("H3").Value = 4500
$ws.Range("I3").Value = "W"

$ws.Range("B3").Select()
